$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters used: D=4, L=12, M=13, N=14, O=15, P=16, Q=17, S=19, T=20
$rows = @(
  [pscustomobject]@{Row=23; D=44635; L='Primera'; M=40; N=23000; O=23000; P=23000; Q='$/bandeja 18 kilos granel'; S=1278; T=18},
  [pscustomobject]@{Row=24; D=44414; L='Primera'; M=95; N=12000; O=13000; P=12526; Q='$/bandeja 18 kilos granel'; S=696; T=18},
  [pscustomobject]@{Row=25; D=44293; L='Primera'; M=50; N=13000; O=13000; P=13000; Q='$/bandeja 18 kilos granel'; S=722; T=18},
  [pscustomobject]@{Row=26; D=44309; L='Primera'; M=80; N=13000; O=13000; P=13000; Q='$/bandeja 18 kilos granel'; S=722; T=18},
  [pscustomobject]@{Row=27; D=44266; L='Primera'; M=100; N=15000; O=15000; P=15000; Q='$/bandeja 18 kilos granel'; S=833; T=18},
  [pscustomobject]@{Row=28; D=44286; L='Primera'; M=65; N=13000; O=13000; P=13000; Q='$/bandeja 18 kilos granel'; S=722; T=18},
  [pscustomobject]@{Row=29; D=44379; L='Primera'; M=80; N=12000; O=13000; P=12500; Q='$/bandeja 18 kilos granel'; S=694; T=18},
  [pscustomobject]@{Row=30; D=44259; L='Primera'; M=150; N=14000; O=15000; P=14567; Q='$/bandeja 18 kilos granel'; S=809; T=18},
  [pscustomobject]@{Row=31; D=44620; L='Primera'; M=5; N=360000; O=360000; P=360000; Q='$/bins (450 kilos)'; S=800; T=450},
  [pscustomobject]@{Row=32; D=44620; L='Primera'; M=65; N=16000; O=16000; P=16000; Q='$/caja 18 kilos granel'; S=889; T=18},
  [pscustomobject]@{Row=33; D=44308; L='Primera'; M=150; N=13000; O=15000; P=14067; Q='$/bandeja 18 kilos granel'; S=782; T=18},
  [pscustomobject]@{Row=34; D=44308; L='Primera'; M=90; N=17000; O=17000; P=17000; Q='$/caja 18 kilos empedrada'; S=944; T=18},
  [pscustomobject]@{Row=35; D=44427; L='Primera'; M=80; N=13000; O=13000; P=13000; Q='$/bandeja 18 kilos granel'; S=722; T=18},
  [pscustomobject]@{Row=36; D=44623; L='Primera'; M=90; N=17000; O=18000; P=17389; Q='$/bandeja 18 kilos granel'; S=966; T=18},
  [pscustomobject]@{Row=37; D=44453; L='Especial'; M=55; N=20000; O=20000; P=20000; Q='$/bandeja 18 kilos granel'; S=1111; T=18},
  [pscustomobject]@{Row=38; D=44624; L='Primera'; M=55; N=18000; O=18000; P=18000; Q='$/bandeja 18 kilos granel'; S=1000; T=18},
  [pscustomobject]@{Row=39; D=44340; L='Primera'; M=105; N=12000; O=13000; P=12619; Q='$/bandeja 18 kilos granel'; S=701; T=18},
  [pscustomobject]@{Row=40; D=44434; L='Primera'; M=40; N=13000; O=13000; P=13000; Q='$/bandeja 18 kilos granel'; S=722; T=18},
  [pscustomobject]@{Row=41; D=44270; L='Primera'; M=65; N=15000; O=15000; P=15000; Q='$/bandeja 18 kilos granel'; S=833; T=18},
  [pscustomobject]@{Row=42; D=44270; L='Primera'; M=5; N=250000; O=250000; P=250000; Q='$/bins (450 kilos)'; S=556; T=450},
  [pscustomobject]@{Row=43; D=44391; L='Primera'; M=80; N=12000; O=13000; P=12500; Q='$/bandeja 18 kilos granel'; S=694; T=18},
  [pscustomobject]@{Row=44; D=44376; L='Primera'; M=95; N=13000; O=13000; P=13000; Q='$/bandeja 18 kilos granel'; S=722; T=18},
  [pscustomobject]@{Row=45; D=44307; L='Primera'; M=80; N=13000; O=13000; P=13000; Q='$/bandeja 18 kilos granel'; S=722; T=18},
  [pscustomobject]@{Row=46; D=44622; L='Primera'; M=55; N=16000; O=16000; P=16000; Q='$/bandeja 18 kilos granel'; S=889; T=18},
  [pscustomobject]@{Row=47; D=44622; L='Primera'; M=3; N=360000; O=360000; P=360000; Q='$/bins (450 kilos)'; S=800; T=450},
  [pscustomobject]@{Row=48; D=44405; L='Primera'; M=50; N=13000; O=13000; P=13000; Q='$/bandeja 18 kilos granel'; S=722; T=18},
  [pscustomobject]@{Row=49; D=44377; L='Primera'; M=80; N=12000; O=13000; P=12500; Q='$/bandeja 18 kilos granel'; S=694; T=18},
  [pscustomobject]@{Row=50; D=44424; L='Primera'; M=45; N=13000; O=13000; P=13000; Q='$/bandeja 18 kilos granel'; S=722; T=18},
  [pscustomobject]@{Row=51; D=44306; L='Especial'; M=40; N=18000; O=18000; P=18000; Q='$/bandeja 18 kilos granel'; S=1000; T=18},
  [pscustomobject]@{Row=52; D=44384; L='Primera'; M=95; N=12000; O=12000; P=12000; Q='$/bandeja 18 kilos granel'; S=667; T=18},
  [pscustomobject]@{Row=53; D=44383; L='Primera'; M=120; N=12000; O=13000; P=12542; Q='$/bandeja 18 kilos granel'; S=697; T=18},
  [pscustomobject]@{Row=54; D=44295; L='Primera'; M=40; N=13000; O=13000; P=13000; Q='$/bandeja 18 kilos granel'; S=722; T=18},
  [pscustomobject]@{Row=55; D=44369; L='Primera'; M=105; N=12000; O=13000; P=12381; Q='$/bandeja 18 kilos granel'; S=688; T=18},
  [pscustomobject]@{Row=56; D=44369; L='Segunda'; M=65; N=8000; O=8000; P=8000; Q='$/bandeja 18 kilos granel'; S=444; T=18},
  [pscustomobject]@{Row=57; D=44278; L='Primera'; M=80; N=14000; O=15000; P=14500; Q='$/bandeja 18 kilos granel'; S=806; T=18},
  [pscustomobject]@{Row=58; D=44435; L='Primera'; M=140; N=13000; O=13000; P=13000; Q='$/bandeja 18 kilos granel'; S=722; T=18},
  [pscustomobject]@{Row=59; D=44292; L='Primera'; M=50; N=14000; O=14000; P=14000; Q='$/bandeja 18 kilos granel'; S=778; T=18},
  [pscustomobject]@{Row=60; D=44314; L='Especial'; M=35; N=20000; O=20000; P=20000; Q='$/bandeja 18 kilos granel'; S=1111; T=18},
  [pscustomobject]@{Row=61; D=44314; L='Primera'; M=55; N=14000; O=14000; P=14000; Q='$/bandeja 18 kilos granel'; S=778; T=18},
  [pscustomobject]@{Row=62; D=44354; L='Primera'; M=120; N=13000; O=14000; P=13542; Q='$/bandeja 18 kilos granel'; S=752; T=18},
  [pscustomobject]@{Row=63; D=44354; L='Primera'; M=5; N=270000; O=270000; P=270000; Q='$/bins (450 kilos)'; S=600; T=450},
  [pscustomobject]@{Row=64; D=44406; L='Primera'; M=70; N=12000; O=13000; P=12429; Q='$/bandeja 18 kilos granel'; S=690; T=18},
  [pscustomobject]@{Row=65; D=44260; L='Primera'; M=65; N=15000; O=15000; P=15000; Q='$/bandeja 18 kilos granel'; S=833; T=18},
  [pscustomobject]@{Row=66; D=44341; L='Especial'; M=65; N=15000; O=15000; P=15000; Q='$/bandeja 18 kilos granel'; S=833; T=18},
  [pscustomobject]@{Row=67; D=44341; L='Primera'; M=80; N=12000; O=12000; P=12000; Q='$/bandeja 18 kilos granel'; S=667; T=18},
  [pscustomobject]@{Row=68; D=44245; L='Primera'; M=65; N=17000; O=17000; P=17000; Q='$/bandeja 18 kilos granel'; S=944; T=18},
  [pscustomobject]@{Row=69; D=44305; L='Primera'; M=80; N=13000; O=13000; P=13000; Q='$/bandeja 18 kilos granel'; S=722; T=18},
  [pscustomobject]@{Row=70; D=44442; L='Primera'; M=75; N=12000; O=12000; P=12000; Q='$/bandeja 18 kilos granel'; S=667; T=18},
  [pscustomobject]@{Row=71; D=44246; L='Primera'; M=55; N=17000; O=17000; P=17000; Q='$/bandeja 18 kilos granel'; S=944; T=18},
  [pscustomobject]@{Row=72; D=44323; L='Primera'; M=40; N=13000; O=13000; P=13000; Q='$/bandeja 18 kilos granel'; S=722; T=18},
  [pscustomobject]@{Row=73; D=44398; L='Primera'; M=55; N=13000; O=13000; P=13000; Q='$/bandeja 18 kilos granel'; S=722; T=18},
  [pscustomobject]@{Row=74; D=44392; L='Primera'; M=90; N=12000; O=13000; P=12444; Q='$/bandeja 18 kilos granel'; S=691; T=18},
  [pscustomobject]@{Row=75; D=44328; L='Primera'; M=65; N=13000; O=13000; P=13000; Q='$/bandeja 18 kilos granel'; S=722; T=18},
  [pscustomobject]@{Row=76; D=44433; L='Primera'; M=80; N=13000; O=13000; P=13000; Q='$/bandeja 18 kilos granel'; S=722; T=18},
  [pscustomobject]@{Row=77; D=44382; L='Primera'; M=115; N=11000; O=13000; P=12565; Q='$/bandeja 18 kilos granel'; S=698; T=18},
  [pscustomobject]@{Row=78; D=44265; L='Primera'; M=40; N=15000; O=15000; P=15000; Q='$/bandeja 18 kilos granel'; S=833; T=18},
  [pscustomobject]@{Row=79; D=44363; L='Primera'; M=80; N=13000; O=13000; P=13000; Q='$/bandeja 18 kilos granel'; S=722; T=18},
  [pscustomobject]@{Row=80; D=44386; L='Primera'; M=95; N=12000; O=12000; P=12000; Q='$/bandeja 18 kilos granel'; S=667; T=18},
  [pscustomobject]@{Row=81; D=44322; L='Primera'; M=200; N=13000; O=13000; P=13000; Q='$/bandeja 18 kilos granel'; S=722; T=18},
  [pscustomobject]@{Row=82; D=44358; L='Primera'; M=80; N=12000; O=12000; P=12000; Q='$/bandeja 18 kilos granel'; S=667; T=18},
  [pscustomobject]@{Row=83; D=44313; L='Especial'; M=125; N=20000; O=20000; P=20000; Q='$/caja 18 kilos empedrada'; S=1111; T=18},
  [pscustomobject]@{Row=84; D=44244; L='Primera'; M=45; N=17000; O=17000; P=17000; Q='$/bandeja 18 kilos granel'; S=944; T=18},
  [pscustomobject]@{Row=85; D=44356; L='Primera'; M=110; N=12000; O=12000; P=12000; Q='$/bandeja 18 kilos granel'; S=667; T=18},
  [pscustomobject]@{Row=86; D=44302; L='Primera'; M=95; N=11000; O=12000; P=11474; Q='$/bandeja 18 kilos granel'; S=637; T=18},
  [pscustomobject]@{Row=87; D=44291; L='Primera'; M=40; N=13000; O=13000; P=13000; Q='$/bandeja 18 kilos granel'; S=722; T=18},
  [pscustomobject]@{Row=88; D=44326; L='Primera'; M=115; N=13000; O=14000; P=13565; Q='$/bandeja 18 kilos granel'; S=754; T=18},
  [pscustomobject]@{Row=89; D=44348; L='Primera'; M=40; N=13000; O=13000; P=13000; Q='$/bandeja 18 kilos granel'; S=722; T=18},
  [pscustomobject]@{Row=90; D=44281; L='Primera'; M=60; N=13000; O=14000; P=13417; Q='$/bandeja 18 kilos granel'; S=745; T=18},
  [pscustomobject]@{Row=91; D=44271; L='Especial'; M=35; N=20000; O=20000; P=20000; Q='$/bandeja 18 kilos granel'; S=1111; T=18},
  [pscustomobject]@{Row=92; D=44271; L='Primera'; M=95; N=15000; O=15000; P=15000; Q='$/bandeja 18 kilos granel'; S=833; T=18},
  [pscustomobject]@{Row=93; D=44420; L='Primera'; M=65; N=13000; O=13000; P=13000; Q='$/bandeja 18 kilos granel'; S=722; T=18},
  [pscustomobject]@{Row=94; D=44343; L='Primera'; M=75; N=13000; O=13000; P=13000; Q='$/bandeja 18 kilos granel'; S=722; T=18},
  [pscustomobject]@{Row=95; D=44315; L='Primera'; M=85; N=13000; O=14000; P=13529; Q='$/bandeja 18 kilos granel'; S=752; T=18},
  [pscustomobject]@{Row=96; D=44315; L='Segunda'; M=55; N=9000; O=9000; P=9000; Q='$/bandeja 18 kilos granel'; S=500; T=18},
  [pscustomobject]@{Row=97; D=44336; L='Primera'; M=80; N=12000; O=12000; P=12000; Q='$/bandeja 18 kilos granel'; S=667; T=18},
  [pscustomobject]@{Row=98; D=44336; L='Primera'; M=3; N=230000; O=230000; P=230000; Q='$/bins (450 kilos)'; S=511; T=450},
  [pscustomobject]@{Row=99; D=44400; L='Primera'; M=120; N=12000; O=13000; P=12542; Q='$/bandeja 18 kilos granel'; S=697; T=18},
  [pscustomobject]@{Row=100; D=44627; L='Primera'; M=80; N=16000; O=16000; P=16000; Q='$/bandeja 18 kilos granel'; S=889; T=18},
  [pscustomobject]@{Row=101; D=44334; L='Primera'; M=90; N=12000; O=12000; P=12000; Q='$/bandeja 18 kilos granel'; S=667; T=18},
  [pscustomobject]@{Row=102; D=44319; L='Primera'; M=180; N=12000; O=13000; P=12444; Q='$/bandeja 18 kilos granel'; S=691; T=18},
  [pscustomobject]@{Row=103; D=44280; L='Primera'; M=200; N=12000; O=12000; P=12000; Q='$/bandeja 18 kilos granel'; S=667; T=18},
  [pscustomobject]@{Row=104; D=44280; L='Primera'; M=3; N=240000; O=240000; P=240000; Q='$/bins (450 kilos)'; S=533; T=450},
  [pscustomobject]@{Row=105; D=44362; L='Primera'; M=40; N=13000; O=13000; P=13000; Q='$/bandeja 18 kilos granel'; S=722; T=18},
  [pscustomobject]@{Row=106; D=44431; L='Primera'; M=20; N=13000; O=13000; P=13000; Q='$/bandeja 18 kilos granel'; S=722; T=18},
  [pscustomobject]@{Row=107; D=44365; L='Primera'; M=70; N=12000; O=13000; P=12571; Q='$/bandeja 18 kilos granel'; S=698; T=18},
  [pscustomobject]@{Row=108; D=44357; L='Primera'; M=125; N=12000; O=12000; P=12000; Q='$/bandeja 18 kilos granel'; S=667; T=18},
  [pscustomobject]@{Row=109; D=44397; L='Primera'; M=85; N=13000; O=13000; P=13000; Q='$/bandeja 18 kilos granel'; S=722; T=18},
  [pscustomobject]@{Row=110; D=44446; L='Especial'; M=40; N=20000; O=20000; P=20000; Q='$/bandeja 18 kilos granel'; S=1111; T=18},
  [pscustomobject]@{Row=111; D=44329; L='Primera'; M=115; N=13000; O=14000; P=13435; Q='$/bandeja 18 kilos granel'; S=746; T=18},
  [pscustomobject]@{Row=112; D=44637; L='Primera'; M=85; N=23000; O=23000; P=23000; Q='$/bandeja 18 kilos granel'; S=1278; T=18},
  [pscustomobject]@{Row=113; D=44355; L='Primera'; M=115; N=13000; O=14000; P=13435; Q='$/bandeja 18 kilos granel'; S=746; T=18},
  [pscustomobject]@{Row=114; D=44294; L='Primera'; M=100; N=13000; O=13000; P=13000; Q='$/bandeja 18 kilos granel'; S=722; T=18},
  [pscustomobject]@{Row=115; D=44617; L='Primera'; M=20; N=16000; O=16000; P=16000; Q='$/bandeja 18 kilos granel'; S=889; T=18},
  [pscustomobject]@{Row=116; D=44264; L='Primera'; M=30; N=16000; O=16000; P=16000; Q='$/bandeja 18 kilos granel'; S=889; T=18},
  [pscustomobject]@{Row=117; D=44396; L='Primera'; M=175; N=12000; O=13000; P=12457; Q='$/bandeja 18 kilos granel'; S=692; T=18},
  [pscustomobject]@{Row=118; D=44279; L='Primera'; M=120; N=12000; O=13000; P=12417; Q='$/bandeja 18 kilos granel'; S=690; T=18},
  [pscustomobject]@{Row=119; D=44330; L='Primera'; M=95; N=13000; O=14000; P=13526; Q='$/bandeja 18 kilos granel'; S=751; T=18},
  [pscustomobject]@{Row=120; D=44301; L='Primera'; M=45; N=13000; O=13000; P=13000; Q='$/bandeja 18 kilos granel'; S=722; T=18},
  [pscustomobject]@{Row=121; D=44370; L='Primera'; M=65; N=12000; O=12000; P=12000; Q='$/bandeja 18 kilos granel'; S=667; T=18},
  [pscustomobject]@{Row=122; D=44385; L='Primera'; M=285; N=12000; O=13000; P=12561; Q='$/bandeja 18 kilos granel'; S=698; T=18},
  [pscustomobject]@{Row=123; D=44413; L='Primera'; M=95; N=12000; O=13000; P=12579; Q='$/bandeja 18 kilos granel'; S=699; T=18},
  [pscustomobject]@{Row=124; D=44312; L='Primera'; M=210; N=13000; O=14000; P=13262; Q='$/caja 18 kilos granel'; S=737; T=18},
  [pscustomobject]@{Row=125; D=44399; L='Primera'; M=95; N=12000; O=12000; P=12000; Q='$/bandeja 18 kilos granel'; S=667; T=18},
  [pscustomobject]@{Row=126; D=44615; L='Primera'; M=50; N=16000; O=16000; P=16000; Q='$/bandeja 18 kilos granel'; S=889; T=18},
  [pscustomobject]@{Row=127; D=44277; L='Primera'; M=160; N=13000; O=15000; P=14000; Q='$/bandeja 18 kilos granel'; S=778; T=18},
  [pscustomobject]@{Row=128; D=44258; L='Primera'; M=110; N=15000; O=15000; P=15000; Q='$/bandeja 18 kilos granel'; S=833; T=18},
  [pscustomobject]@{Row=129; D=44390; L='Primera'; M=140; N=12000; O=13000; P=12571; Q='$/bandeja 18 kilos granel'; S=698; T=18},
  [pscustomobject]@{Row=130; D=44349; L='Primera'; M=30; N=13000; O=13000; P=13000; Q='$/bandeja 18 kilos granel'; S=722; T=18},
  [pscustomobject]@{Row=131; D=44285; L='Primera'; M=55; N=13000; O=13000; P=13000; Q='$/bandeja 18 kilos granel'; S=722; T=18},
  [pscustomobject]@{Row=132; D=44335; L='Primera'; M=90; N=12000; O=13000; P=12556; Q='$/bandeja 18 kilos granel'; S=698; T=18}
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 12).Value = $r.L
    $ws.Cells.Item($r.Row, 13).Value = $r.M
    $ws.Cells.Item($r.Row, 14).Value = $r.N
    $ws.Cells.Item($r.Row, 15).Value = $r.O
    $ws.Cells.Item($r.Row, 16).Value = $r.P
    $ws.Cells.Item($r.Row, 17).Value = $r.Q
    $ws.Cells.Item($r.Row, 19).Value = $r.S
    $ws.Cells.Item($r.Row, 20).Value = $r.T
}

# Fill the constant (non-shifting) columns for the two brand-new rows (131, 132)
$newRowsConst = @(
  [pscustomobject]@{Row=131; A=10; B='Vega Modelo de Temuco'; C='La Araucanía'; E=9; F='Fruta'; G=100104; H='Frutos de pepita'; I=100104003; J='Membrillo'; K='Champion'; R='Región de O''Higgins'},
  [pscustomobject]@{Row=132; A=10; B='Vega Modelo de Temuco'; C='La Araucanía'; E=9; F='Fruta'; G=100104; H='Frutos de pepita'; I=100104003; J='Membrillo'; K='Champion'; R='Región de O''Higgins'}
)

foreach ($r in $newRowsConst) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
    $ws.Cells.Item($r.Row, 8).Value = $r.H
    $ws.Cells.Item($r.Row, 9).Value = $r.I
    $ws.Cells.Item($r.Row, 10).Value = $r.J
    $ws.Cells.Item($r.Row, 11).Value = $r.K
    $ws.Cells.Item($r.Row, 18).Value = $r.R
}
